# Gacha.xlsx edit: ShopEquipTable / GachaEquipTable update
#  - reorder tabs so ShopEquipTable comes before GachaEquipTable
#  - ShopEquipTable: replace the "Equip50" row with a new "Equip20" row
#  - GachaEquipTable: rebuild the table with a new "rarity" column and new odds
#  - selection / active-tab bookkeeping to match

$wb = $excel.ActiveWorkbook

# --- 1. Reorder worksheets: ShopEquipTable moves in front of GachaEquipTable ---
$shopEquip  = $wb.Worksheets.Item("ShopEquipTable")
$gachaEquip = $wb.Worksheets.Item("GachaEquipTable")
$shopEquip.Move($gachaEquip)

# --- 2. ShopEquipTable: swap the Equip50 row (50 / 800) for an Equip20 row (20 / 300) ---
$ws = $wb.Worksheets.Item("ShopEquipTable")

$ws.Range("A2").Value = "Equip1"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 20

$ws.Range("A3").Value = "Equip10"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 175

$ws.Range("A4").Value = "Equip20"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 300

# --- 3. GachaEquipTable: new grade / rarity / prob table (adds a "rarity" column) ---
$ws2 = $wb.Worksheets.Item("GachaEquipTable")

$ws2.Range("A1").Value = "grade|Int"
$ws2.Range("B1").Value = "rarity|Int"
$ws2.Range("C1").Value = "prob|float"

$ws2.Range("A2").Value = 3
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 0.005

$ws2.Range("A3").Value = 3
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = 0.025

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 0
$ws2.Range("C4").Value = 0.05

$ws2.Range("A5").Value = 2
$ws2.Range("B5").Value = 0
$ws2.Range("C5").Value = 0.2

$ws2.Range("A6").Value = 1
$ws2.Range("B6").Value = 0
$ws2.Range("C6").Value = 0.32

$ws2.Range("A7").Value = 0
$ws2.Range("B7").Value = 0
$ws2.Range("C7").Value = 0.4

$ws2.Range("A2").Select()

# --- 4. ShopSpellTable: selection moves to B4, it is no longer the selected tab ---
$ws3 = $wb.Worksheets.Item("ShopSpellTable")
$ws3.Range("B4").Select()

# --- 5. ShopEquipTable becomes the active tab/sheet (workbook activeTab -> 4) ---
$ws.Range("A2").Select()
